$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7062078129466158
$ws.Range("C2").Value = 0.1110542594521462
$ws.Range("D2").Value = 0.1351963726444581
$ws.Range("F2").Value = 2.936124529211241
$ws.Range("G2").Value = 0.002443244553283323
$ws.Range("M2").Value = 0.8870578730246592
$ws.Range("B3").Value = 0.6438362344940458
$ws.Range("C3").Value = 0.09645981523465252
$ws.Range("D3").Value = 0.1241782874393209
$ws.Range("F3").Value = 2.713953354689124
$ws.Range("G3").Value = 0.002449556922653842
$ws.Range("M3").Value = 0.7895017007456175
$ws.Range("B4").Value = 0.6062791358478705
$ws.Range("C4").Value = 0.087579010562024
$ws.Range("D4").Value = 0.1174106975147424
$ws.Range("F4").Value = 2.578230831821116
$ws.Range("G4").Value = 0.002453625511293309
$ws.Range("M4").Value = 0.7301836773989265
$ws.Range("B5").Value = 0.5911579579166357
$ws.Range("C5").Value = 0.08397943759374016
$ws.Range("D5").Value = 0.1146515441342189
$ws.Range("F5").Value = 2.523086278532332
$ws.Range("G5").Value = 0.002455332174300686
$ws.Range("M5").Value = 0.7061498467686391
$ws.Range("B6").Value = 0.5886581246253115
$ws.Range("C6").Value = 0.08338288293327878
$ws.Range("D6").Value = 0.1141932898567575
$ws.Range("F6").Value = 2.513939139444972
$ws.Range("G6").Value = 0.002455618510408565
$ws.Range("M6").Value = 0.7021672092802618
$ws.Range("B7").Value = 0.6060744657241628
$ws.Range("C7").Value = 0.08753038778468181
$ws.Range("D7").Value = 0.1173734928043899
$ws.Range("F7").Value = 2.577486484841813
$ws.Range("G7").Value = 0.002453648330729755
$ws.Range("M7").Value = 0.7298589968649623
$ws.Range("B8").Value = 0.6845472528653715
$ws.Range("C8").Value = 0.1060049801290575
$ws.Range("D8").Value = 0.1313974228434915
$ws.Range("F8").Value = 2.859370034301747
$ws.Range("G8").Value = 0.002445381174069589
$ws.Range("M8").Value = 0.8532951081362796
$ws.Range("B9").Value = 0.8444113891143274
$ws.Range("C9").Value = 0.1429065185907348
$ws.Range("D9").Value = 0.1589113843968022
$ws.Range("F9").Value = 3.418128215519602
$ws.Range("G9").Value = 0.002430689355263423
$ws.Range("M9").Value = 1.100329235156664
$ws.Range("B10").Value = 0.9656818871740143
$ws.Range("C10").Value = 0.1704819173461658
$ws.Range("D10").Value = 0.179180617294179
$ws.Range("F10").Value = 3.833030792646525
$ws.Range("G10").Value = 0.00242080858791504
$ws.Range("M10").Value = 1.285383672960236
$ws.Range("B11").Value = 1.021719550366186
$ws.Range("C11").Value = 0.183139461686892
$ws.Range("D11").Value = 0.1884238169948276
$ws.Range("F11").Value = 4.02289593659998
$ws.Range("G11").Value = 0.002416509007142571
$ws.Range("M11").Value = 1.370460808347318
$ws.Range("B12").Value = 1.043067865889441
$ws.Range("C12").Value = 0.1879498017400181
$ws.Range("D12").Value = 0.1919280485178092
$ws.Range("F12").Value = 4.094967313273514
$ws.Range("G12").Value = 0.002414908717598985
$ws.Range("M12").Value = 1.40281566585135
$ws.Range("B13").Value = 1.038464382721486
$ws.Range("C13").Value = 0.1869130292356829
$ws.Range("D13").Value = 0.1911731576870466
$ws.Range("F13").Value = 4.079437540537583
$ws.Range("G13").Value = 0.0024152521323529
$ws.Range("M13").Value = 1.395841179998285
$ws.Range("B14").Value = 1.023473306080746
$ws.Range("C14").Value = 0.1835348612398207
$ws.Range("D14").Value = 0.1887120269085472
$ws.Range("F14").Value = 4.028821744292486
$ws.Range("G14").Value = 0.002416376792893005
$ws.Range("M14").Value = 1.373119842065933
$ws.Range("B15").Value = 1.014307601777944
$ws.Range("C15").Value = 0.1814679032707716
$ws.Range("D15").Value = 0.1872050617155594
$ws.Range("F15").Value = 3.997841069888125
$ws.Range("G15").Value = 0.002417069303383572
$ws.Range("M15").Value = 1.359220618176934
$ws.Range("B16").Value = 0.9620374608120983
$ws.Range("C16").Value = 0.1696570746811972
$ws.Range("D16").Value = 0.1785770740327166
$ws.Range("F16").Value = 3.820646291062104
$ws.Range("G16").Value = 0.002421093486410788
$ws.Range("M16").Value = 1.279842479262058
$ws.Range("B17").Value = 0.9301965025731533
$ws.Range("C17").Value = 0.162441181998588
$ws.Range("D17").Value = 0.1732904281550987
$ws.Range("F17").Value = 3.712239320984139
$ws.Range("G17").Value = 0.00242361204952703
$ws.Range("M17").Value = 1.231382307525422
$ws.Range("B18").Value = 0.9119642452985204
$ws.Range("C18").Value = 0.1583014057489152
$ws.Range("D18").Value = 0.1702517714401068
$ws.Range("F18").Value = 3.649991555529766
$ws.Range("G18").Value = 0.002425079048772746
$ws.Range("M18").Value = 1.203592961534923
$ws.Range("B19").Value = 0.9058050847783647
$ws.Range("C19").Value = 0.1569015445591333
$ws.Range("D19").Value = 0.1692232685646502
$ws.Range("F19").Value = 3.628933214719012
$ws.Range("G19").Value = 0.002425578914018428
$ws.Range("M19").Value = 1.194198056923582
$ws.Range("B20").Value = 0.9335775388524894
$ws.Range("C20").Value = 0.1632082210843464
$ws.Range("D20").Value = 0.1738529795320858
$ws.Range("F20").Value = 3.723768469192123
$ws.Range("G20").Value = 0.002423342042569781
$ws.Range("M20").Value = 1.236532250248416
$ws.Range("B21").Value = 1.027873050161588
$ws.Range("C21").Value = 0.1845266368880232
$ws.Range("D21").Value = 0.1894348049698351
$ws.Range("F21").Value = 4.043684023620244
$ws.Range("G21").Value = 0.002416045698390501
$ws.Range("M21").Value = 1.379789829170406
$ws.Range("B22").Value = 1.090248462009754
$ws.Range("C22").Value = 0.1985601598350968
$ws.Range("D22").Value = 0.1996423029825394
$ws.Range("F22").Value = 4.25378472819483
$ws.Range("G22").Value = 0.002411439457174026
$ws.Range("M22").Value = 1.474225788383436
$ws.Range("B23").Value = 1.056888154274191
$ws.Range("C23").Value = 0.1910606858326958
$ws.Range("D23").Value = 0.1941919240595666
$ws.Range("F23").Value = 4.141552926314205
$ws.Range("G23").Value = 0.002413883108412396
$ws.Range("M23").Value = 1.423746381884953
$ws.Range("B24").Value = 0.932048743520113
$ws.Range("C24").Value = 0.1628614155690116
$ws.Range("D24").Value = 0.1735986479126268
$ws.Range("F24").Value = 3.718555902831326
$ws.Range("G24").Value = 0.00242346405354921
$ws.Range("M24").Value = 1.23420374053147
$ws.Range("B25").Value = 0.8005061617689648
$ws.Range("C25").Value = 0.1328461689817004
$ws.Range("D25").Value = 0.1514619330618387
$ws.Range("F25").Value = 3.266253382406745
$ws.Range("G25").Value = 0.002434502528124131
$ws.Range("M25").Value = 1.032913852734069
